$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.691.84"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "1.894.66"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.98"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4866"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3795"
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07337"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9162"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07695"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").Value = "1.858.96"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.475"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.618"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.00"

$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008803"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").Value = "27.728.28"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.123"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "2.126.64"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.907"
$ws.Range("E25").Value = "  -2.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.65"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.35"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.154"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.72"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.906"
$ws.Range("E30").Value = "  -1.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08916"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  -5.36%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7654"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.646"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02043"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.523"
$ws.Range("E38").Value = "  -7.24%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.093"
$ws.Range("E39").Value = "  -3.56%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05270"
$ws.Range("E40").Value = "  -1.81%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5467"
$ws.Range("E41").Value = "  -3.32%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.986"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.911"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1517"
$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.443"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "110.08"
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.63"
$ws.Range("E47").Value = "  -0.55%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4790"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.637"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.35"
$ws.Range("E51").Value = "  -0.62%  "

